$d = $word.ActiveDocument

# --- Highlight "Ecoregion" (without the trailing "s") in yellow ---
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$rng.Find.Replacement.Highlight = $true
$rng.Find.Execute("Ecoregion", $false, $false, $false, $false, $false, $true, 1, $true, "Ecoregion", 2) | Out-Null

# --- Highlight "(mixed effects model" (without the trailing ")." ) in yellow ---
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Replacement.ClearFormatting()
$rng2.Find.Replacement.Highlight = $true
$rng2.Find.Execute("(mixed effects model", $false, $false, $false, $false, $false, $true, 1, $true, "(mixed effects model", 2) | Out-Null

# --- Move the "_GoBack" bookmark from the end of the document to right after
#     "Ecoregion" (splitting "Ecoregion" and the following "s") ---
$old = $d.Bookmarks("_GoBack")
$old.Delete()

$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$found = $rng3.Find.Execute("Ecoregion", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $d.Range($rng3.End, $rng3.End)) | Out-Null

$d.Save()
